# Update Work Week and Social Spending
# (Refresh the GDP per Capita data series for Ukraine with updated values
# and extend the series through 2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# All of the "Data" column values in this sheet are stored as text (not
# numbers), so force the text number format before writing the values,
# then clear the formatting back off again so the cells keep the
# workbook's default (unstyled) appearance.
$dataRange = $ws.Range("E2:E45")
$dataRange.NumberFormat = "@"

# Row 2 (year 1973): updated data value
$ws.Range("E2").Value = "7849"

# Rows 9-39 (years 1980-2010): refreshed data values
$ws.Range("E9").Value = "8467"
$ws.Range("E10").Value = "8695"
$ws.Range("E11").Value = "8934"
$ws.Range("E12").Value = "9288"
$ws.Range("E13").Value = "9476"
$ws.Range("E14").Value = "9387"
$ws.Range("E15").Value = "9408"
$ws.Range("E16").Value = "9648"
$ws.Range("E17").Value = "9704"
$ws.Range("E18").Value = "10082"
$ws.Range("E19").Value = "9607"
$ws.Range("E20").Value = "8899.09669056187"
$ws.Range("E21").Value = "8130.2767858566"
$ws.Range("E22").Value = "7090.25734143651"
$ws.Range("E23").Value = "5585.41332198664"
$ws.Range("E24").Value = "5024.18173524013"
$ws.Range("E25").Value = "4636.62630920416"
$ws.Range("E26").Value = "4612.40421983972"
$ws.Range("E27").Value = "4640.19919778576"
$ws.Range("E28").Value = "4751.31624059898"
$ws.Range("E29").Value = "5165.80815372985"
$ws.Range("E30").Value = "5793.65975965037"
$ws.Range("E31").Value = "6262.98242480859"
$ws.Range("E32").Value = "7029.57052379711"
$ws.Range("E33").Value = "8049.39537703245"
$ws.Range("E34").Value = "8496.67633895713"
$ws.Range("E35").Value = "9355.98792545806"
$ws.Range("E36").Value = "10355.8184466736"
$ws.Range("E37").Value = "10824.07646107"
$ws.Range("E38").Value = "9381.2249754413"
$ws.Range("E39").Value = "9601.20860859582"

# New rows 40-45 (years 2011-2016): extend the series
$ws.Range("A40").Value = 804
$ws.Range("B40").Value = "Ukraine"
$ws.Range("C40").Value = "GDP per Capita"
$ws.Range("D40").Value = 2011
$ws.Range("E40").Value = "10333"

$ws.Range("A41").Value = 804
$ws.Range("B41").Value = "Ukraine"
$ws.Range("C41").Value = "GDP per Capita"
$ws.Range("D41").Value = 2012
$ws.Range("E41").Value = "10383"

$ws.Range("A42").Value = 804
$ws.Range("B42").Value = "Ukraine"
$ws.Range("C42").Value = "GDP per Capita"
$ws.Range("D42").Value = 2013
$ws.Range("E42").Value = "10404"

$ws.Range("A43").Value = 804
$ws.Range("B43").Value = "Ukraine"
$ws.Range("C43").Value = "GDP per Capita"
$ws.Range("D43").Value = 2014
$ws.Range("E43").Value = "9818"

$ws.Range("A44").Value = 804
$ws.Range("B44").Value = "Ukraine"
$ws.Range("C44").Value = "GDP per Capita"
$ws.Range("D44").Value = 2015
$ws.Range("E44").Value = "8961"

$ws.Range("A45").Value = 804
$ws.Range("B45").Value = "Ukraine"
$ws.Range("C45").Value = "GDP per Capita"
$ws.Range("D45").Value = 2016
$ws.Range("E45").Value = "9214"

# Restore the default (unstyled) formatting on the column now that the
# text values have been written.
$dataRange.ClearFormats()
